# Daily attendance processing - 2025-11-16 14:18:46
#
# Rotates the "Recorded By" author list (column G) so that the literal
# "System" entry that currently sits at the front of the comma-separated
# list is moved to the end (equivalent to a right-rotation of the list).
# Only cells whose text matches one of the known before-states are touched,
# which mirrors exactly which rows changed upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of exact current value -> new value for the "Recorded By" column (G)
$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
